$d = $word.ActiveDocument

# 1. Remove the whole "Meta description: ..." paragraph (the second paragraph
#    of the document, right after the Heading1 title paragraph). Deleting the
#    paragraph's full range (text + its trailing paragraph mark) collapses the
#    <w:p> cleanly, merging the following paragraph up into its place.
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range
$metaRange.Delete()

# 2. Split the final paragraph (the "Please create a feature image..." image
#    prompt paragraph) into two paragraphs:
#      - a new bold paragraph with the page title text
#      - the existing italic paragraph, but with its text replaced by the
#        former meta-description copy
#    Replacing the whole paragraph (including its paragraph mark) with a
#    two-paragraph WordOpenXML fragment lets us control the run layout
#    exactly (leading empty run + formatted run) instead of inheriting
#    formatting from a split point.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xmlFragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Deadly 5 Free - Wild West Themed Slot Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Deadly 5, a wild west-themed slot game featuring four outlaw characters. Play for free and experience the excitement of the American frontier.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastRange.InsertXML($xmlFragment)
